$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("G2").Value = 88
    $ws.Range("G3").Value = 78
    $ws.Range("F4").Value = 8643
    $ws.Range("G4").Value = 78
    $ws.Range("G5").Value = "不可售"
    $ws.Range("G6").Value = 88
}

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F14").Value = 76
$ws1.Range("F17").Value = 6109
$ws1.Range("F18").Value = 204
$ws1.Range("F19").Value = 302
$ws1.Range("F20").Value = 2244
$ws1.Range("F21").Value = 88
$ws1.Range("F22").Value = 151

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F16").Value = 76
$ws4.Range("F20").Value = 6109
$ws4.Range("F22").Value = 204
$ws4.Range("F23").Value = 302
$ws4.Range("F24").Value = 2244
$ws4.Range("F25").Value = 88
$ws4.Range("F26").Value = 151
